$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new "MVC - Core Sports Store App, 1 - Chapter 9" row was added to the
# Work Log table, right after the existing "...Chapter 8" row (row 24) and
# before "Assignment Q&A Week 5" (old row 25). Inserting a whole sheet row
# there shifts the three rows below it (Assignment Q&A Week 5, Labs 4, ADA
# Compliance) down by one, from rows 25-27 to rows 26-28.
$ws.Rows.Item(25).Insert()

# Populate the new row 25 with the new entry (matches the formatting of the
# other "MVC - Core Sports Store App, 1 - Chapter N" rows above it: bold
# "s=5" style on the name cell, centered style on hours, centered/date style
# on the due date).
$ws.Range("B25").Value = "MVC - Core Sports Store App, 1 - Chapter 9"
$ws.Range("C25").Value = 4
# The "Due Date" column for this block of rows uses a leading non-breaking
# space before the date text (matches the existing shared string used by the
# surrounding rows instead of creating a new near-duplicate string).
$ws.Range("D25").Value = ([char]0x00A0 + "5/10/2020")

# The commit also fills in the previously-blank "Time in hrs" values for the
# "Labs 4" and "ADA Compliance" rows (now pushed down to rows 27 and 28).
$ws.Range("C27").Value = 1
$ws.Range("C28").Value = 1

# Grow Table1 so the new row participates in the table / autofilter range.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:F28"))

# Leave the cursor where the author ended up after the edit.
$ws.Range("D30").Select() | Out-Null
